$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 116; existing rows 116-154 shift down to 117-155.
$ws.Rows.Item(116).Insert()

# Populate the newly inserted row 116 with the new weekly price record.
$ws.Range("A116").Value = 11
$ws.Range("B116").Value = "Vega Monumental Concepción"
$ws.Range("C116").Value = "Bíobío"
$ws.Range("D116").Value = 44468
$ws.Range("E116").Value = 8
$ws.Range("F116").Value = "Fruta"
$ws.Range("G116").Value = 100101
$ws.Range("H116").Value = "Berries"
$ws.Range("I116").Value = 100112025
$ws.Range("J116").Value = "Frutilla"
$ws.Range("K116").Value = "Sin especificar"
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 100
$ws.Range("N116").Value = 16000
$ws.Range("O116").Value = 16000
$ws.Range("P116").Value = 16000
$ws.Range("Q116").Value = "$/bandeja 7 kilos"
$ws.Range("R116").Value = "Provincia de Melipilla"
$ws.Range("S116").Value = 2286
$ws.Range("T116").Value = 7
